# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 57 (pushing the existing rows 57..156
# down to 58..157) and populate the newly inserted row with the new
# observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 57; this shifts rows 57:156 down to
# 58:157 (carrying their values/formats with them), leaving row 57 empty
# except for the date-format style copied down from the row above.
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with the new record.
$ws.Range("A57").Value = 8
$ws.Range("B57").Value = "Terminal La Palmera de La Serena"
$ws.Range("C57").Value = "Coquimbo"
$ws.Range("D57").Value = 44469
$ws.Range("E57").Value = 4
$ws.Range("F57").Value = 100112012
$ws.Range("G57").Value = "Espinaca"
$ws.Range("H57").Value = "Sin especificar"
$ws.Range("I57").Value = "Primera"
$ws.Range("J57").Value = 2860
$ws.Range("K57").Value = 400
$ws.Range("L57").Value = 500
$ws.Range("M57").Value = 450
$ws.Range("N57").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O57").Value = "Provincia del Elquí"
$ws.Range("P57").Value = 900
$ws.Range("Q57").Value = 0.5
$ws.Range("R57").Value = "Hortaliza"
